$d = $word.ActiveDocument

# The target paragraph holds the "transcription line" id for this page,
# split across three runs: "<id>", "p027r_1", "</id>" (plus a trailing
# empty run). We need to merge the three text runs into a single run
# (keeping the first run's formatting) so the paragraph text becomes the
# single, contiguous string "<id>p027r_1</id>".
$targetText = "<id>p027r_1</id>"
$openTag    = "<id>"

$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    # Paragraph.Range.Text includes the trailing paragraph-mark
    # character(s); strip those before comparing.
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $targetText) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # End of the first run ("<id>") -- everything from here up to the
    # end of the visible text ("</id>") belongs to runs 2 and 3, which
    # get collapsed into run 1.
    $run1End = $r.Start + $openTag.Length
    $textEnd = $r.Start + $targetText.Length

    $toMerge = $d.Range($run1End, $textEnd)
    $toMerge.Delete()

    $run1 = $d.Range($r.Start, $run1End)
    $run1.InsertAfter($targetText.Substring($openTag.Length))
}
